$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (first sheet)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a new row before the existing row 293 ("EQUISAB S.A." under
# "OFICINA-CATAECSA"), which shifts every row below it down by one.
$ws1.Rows.Item(293).Insert()

# Populate the freshly inserted row with the new client record.
$ws1.Range("A293").Value = "OFICINA-CATAECSA"
$ws1.Range("B293").Value = "ECUAFERRI S.A."
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(293, $c).Value = 0
}

# Update the specific sales figures that changed for existing clients.
$ws1.Range("D45").Value = 960.96
$ws1.Range("K69").Value = -406.91
$ws1.Range("M97").Value = -171.19
$ws1.Range("D111").Value = 915.84
$ws1.Range("M111").Value = 1205.63
$ws1.Range("D169").Value = -152.64
$ws1.Range("L169").Value = -1151.4

# The trailing "coverage" summary row moved from 349 to 350 because of the
# inserted row; refresh its "X de 347" -> "X de 348" labels.
$ws1.Range("C350").Value = "0 de 348"
$ws1.Range("D350").Value = "3 de 348"
$ws1.Range("E350").Value = "1 de 348"
$ws1.Range("F350").Value = "0 de 348"
$ws1.Range("G350").Value = "0 de 348"
$ws1.Range("H350").Value = "1 de 348"
$ws1.Range("I350").Value = "1 de 348"
$ws1.Range("J350").Value = "0 de 348"
$ws1.Range("K350").Value = "1 de 348"
$ws1.Range("L350").Value = "3 de 348"
$ws1.Range("M350").Value = "5 de 348"
$ws1.Range("N350").Value = "0 de 348"
$ws1.Range("O350").Value = "0 de 348"
$ws1.Range("P350").Value = "0 de 348"
$ws1.Range("Q350").Value = "0 de 348"
$ws1.Range("R350").Value = "0 de 348"

# ---------------------------------------------------------------------------
# Sheet "VENTA MENSUAL" (second sheet)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Same insertion: a new "ECUAFERRI S.A." row under "OFICINA-CATAECSA"
# before the existing row 297 ("EQUISAB S.A.").
$ws2.Rows.Item(297).Insert()

$ws2.Range("A297").Value = "OFICINA-CATAECSA"
$ws2.Range("B297").Value = "ECUAFERRI S.A."
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(297, $c).Value = 0
}

# Monthly sales figure updates for existing clients.
$ws2.Range("F45").Value = 960.96
$ws2.Range("F69").Value = -406.91
$ws2.Range("F97").Value = 304.01
$ws2.Range("F111").Value = 2121.47
$ws2.Range("F173").Value = -1304.04

# Grand-total row moved from 353 to 354; only the "noviembre" (F) total
# changes to reflect the figure updates above.
$ws2.Range("F354").Value = 11340.32
